$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 12.47849999999999
$ws.Range("E6").Value = 12.5252
$ws.Range("D7").Value = -7.424999999999994
$ws.Range("A10").Value = -20.47329999999997
$ws.Range("A12").Value = -22.72670000000004
$ws.Range("C13").Value = -12.77659999999999
$ws.Range("A18").Value = -22.62260000000003
$ws.Range("D20").Value = -8.685499999999996
